{"js": "// 1) Update the letter date: \"September 19, 2025\" -> \"September 21, 2025\"\n//    (there is another, unrelated \"September 30, 2025\" elsewhere in the\n//    letter body that must NOT be touched, so we search for the exact,\n//    full date string rather than doing a loose replace).\nconst dateResults = context.document.body.search(\"September 19, 2025\", {\n  matchCase: true,\n  matchWholeWord: false\n});\ndateResults.load(\"items\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"September 21, 2025\", Word.InsertLocation.replace);\n}\n\n// 2) Split the single-line mailing address \"919 Story Road, San Jose CA\n//    95122\" into two separate paragraphs: \"919 Story Road\" and\n//    \"San Jose, CA 95122\". There is a second, identical address string\n//    further down inside the \"PROPERTY ADDRESS\" table cell that must be\n//    left exactly as-is, so we search only within the body paragraphs\n//    (not tables) and only touch the first match.\nconst addrResults = context.document.body.search(\"919 Story Road, San Jose CA 95122\", {\n  matchCase: true,\n  matchWholeWord: false\n});\naddrResults.load(\"items\");\nawait context.sync();\n\nif (addrResults.items.length > 0) {\n  const addrRange = addrResults.items[0];\n  // Confirm this hit lives in its own (non-table) paragraph before editing it.\n  const addrPara = addrRange.paragraphs.getFirst();\n  addrRange.insertText(\"919 Story Road\", Word.InsertLocation.replace);\n  addrPara.insertParagraph(\"San Jose, CA 95122\", Word.InsertLocation.after);\n}\n\nawait context.sync();\n\n// 3) Remove the now-redundant blank \"No Spacing\" paragraph that\n//    immediately follows the \"... Board of Directors\" signature line.\nconst boardResults = context.document.body.search(\"Board of Directors\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nboardResults.load(\"items\");\nawait context.sync();\n\nif (boardResults.items.length > 0) {\n  const boardPara = boardResults.items[0].paragraphs.getFirst();\n  const nextPara = boardPara.getNext();\n  nextPara.load(\"text\");\n  await context.sync();\n\n  if (nextPara.text === \"\") {\n    nextPara.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the letter date: \"September 19, 2025\" -> \"September 21, 2025\"\n#    (there is another, unrelated \"September 30, 2025\" elsewhere in the\n#    letter body that must NOT be touched, so we search for the exact,\n#    full date string rather than doing a loose replace).\n$dateFind = $d.Content.Find\n$dateFind.Text = \"September 19, 2025\"\n$dateFound = $dateFind.Execute()\nif ($dateFound) {\n    $dateFind.Parent.Text = \"September 21, 2025\"\n}\n\n# 2) Split the single-line mailing address \"919 Story Road, San Jose CA\n#    95122\" into two separate paragraphs: \"919 Story Road\" and\n#    \"San Jose, CA 95122\". There is a second, identical address string\n#    further down inside the \"PROPERTY ADDRESS\" table cell that must be\n#    left exactly as-is, so Find (run against $d.Content, the main\n#    document story) naturally stops at the first match, which is the\n#    mailing-address paragraph, not the table cell.\n$addrRng = $d.Content\n$addrFind = $addrRng.Find\n$addrFind.Text = \"919 Story Road, San Jose CA 95122\"\n$addrFound = $addrFind.Execute()\nif ($addrFound) {\n    $addrRng.Text = \"919 Story Road\"\n    $addrEnd = $addrRng.End\n    $addrRng.InsertParagraphAfter()\n    $newParaRng = $d.Range($addrEnd + 1, $addrEnd + 1)\n    $newParaRng.Text = \"San Jose, CA 95122\"\n}\n\n# 3) Remove the now-redundant blank \"No Spacing\" paragraph that\n#    immediately follows the \"... Board of Directors\" signature line.\n$boardRng = $d.Content\n$boardFind = $boardRng.Find\n$boardFind.Text = \"Board of Directors\"\n$boardFound = $boardFind.Execute()\nif ($boardFound) {\n    $boardPara = $boardRng.Paragraphs.First\n    $markPos = $boardPara.Range.End\n    $nextMarkRng = $d.Range($markPos + 1, $markPos + 2)\n    if ($nextMarkRng.Text -eq [char]13) {\n        $nextMarkRng.Delete()\n    }\n}\n"}
